$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 326.5
$ws.Range("I2").Value = 218.83333
$ws.Range("J2").Value = 649.5
$ws.Range("K2").Value = 218.83333
$ws.Range("L2").Value = 649.5
$ws.Range("M2").Value = -105.83333
$ws.Range("N2").Value = -875.5
$ws.Range("H12").Value = 1868.5385
$ws.Range("J12").Value = 2299.1
$ws.Range("L12").Value = 2299.1
$ws.Range("N12").Value = -2639.1
$ws.Range("H111").Value = 3866.3333
$ws.Range("I111").Value = 799
$ws.Range("J111").Value = 5400
$ws.Range("K111").Value = 2397
$ws.Range("L111").Value = 16200
$ws.Range("M111").Value = 670
$ws.Range("N111").Value = -22334
$ws.Range("H127").Value = 714.8
$ws.Range("I127").Value = 714.8
$ws.Range("K127").Value = 2144.4
$ws.Range("M127").Value = 2815.6
$ws.Range("H135").Value = 946.1818
$ws.Range("I135").Value = 927.44446
$ws.Range("J135").Value = 1030.5
$ws.Range("K135").Value = 8347.00014
$ws.Range("L135").Value = 9274.5
$ws.Range("M135").Value = -5812.00014
$ws.Range("N135").Value = -14344.5
$ws.Range("H136").Value = 69999
$ws.Range("J136").Value = 69999
$ws.Range("L136").Value = 69999
$ws.Range("N136").Value = -80199
$ws.Range("H138").Value = 3546.3022
$ws.Range("I138").Value = 3024.875
$ws.Range("K138").Value = 9074.625
$ws.Range("M138").Value = -3934.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H139").Value = 69999
$ws.Range("J139").Value = 69999
$ws.Range("L139").Value = 69999
$ws.Range("N139").Value = -80279

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1744.25
$ws.Range("I99").Value = 1743.1
$ws.Range("J99").Value = 1750
$ws.Range("K99").Value = 1743.1
$ws.Range("L99").Value = 1750
$ws.Range("M99").Value = -245.0999999999999
$ws.Range("N99").Value = -4746
$ws.Range("H134").Value = 4329.273
$ws.Range("I134").Value = 1916.5
$ws.Range("J134").Value = 7224.6
$ws.Range("K134").Value = 5749.5
$ws.Range("L134").Value = 21673.8
$ws.Range("M134").Value = -3214.5
$ws.Range("N134").Value = -26743.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 70436.84
$ws.Range("I16").Value = 101218.445
$ws.Range("K16").Value = 101218.445
$ws.Range("M16").Value = -100931.445
$ws.Range("H31").Value = 56328
$ws.Range("J31").Value = 95272.164
$ws.Range("L31").Value = 95272.164
$ws.Range("N31").Value = -95862.164
$ws.Range("H34").Value = 56328
$ws.Range("J34").Value = 95272.164
$ws.Range("L34").Value = 95272.164
$ws.Range("N34").Value = -95676.164
$ws.Range("H55").Value = 13989.5
$ws.Range("I55").Value = 13989.5
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 13989.5
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = -13674.5
$ws.Range("N55").ClearContents()
$ws.Range("H88").Value = 11159.8
$ws.Range("J88").Value = 11159.8
$ws.Range("L88").Value = 11159.8
$ws.Range("N88").Value = -11971.8
$ws.Range("H91").Value = 11159.8
$ws.Range("J91").Value = 11159.8
$ws.Range("L91").Value = 11159.8
$ws.Range("N91").Value = -13967.8
$ws.Range("H113").Value = 70436.84
$ws.Range("I113").Value = 101218.445
$ws.Range("K113").Value = 101218.445
$ws.Range("M113").Value = -99048.44500000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 3709.8667
$ws.Range("I3").Value = 1357.5385
$ws.Range("K3").Value = 4072.6155
$ws.Range("M3").Value = -3960.6155
$ws.Range("H60").Value = 1173974.8
$ws.Range("I60").Value = 1137.1666
$ws.Range("K60").Value = 3411.4998
$ws.Range("M60").Value = -3160.4998
$ws.Range("H104").Value = 5333.3335
$ws.Range("I104").Value = 5000
$ws.Range("J104").Value = 6000
$ws.Range("K104").Value = 15000
$ws.Range("L104").Value = 18000
$ws.Range("M104").Value = -12379
$ws.Range("N104").Value = -23242
$ws.Range("H107").Value = 4112.25
$ws.Range("I107").Value = 3816.3333
$ws.Range("K107").Value = 11448.9999
$ws.Range("M107").Value = -9528.999899999999
$ws.Range("H132").Value = 6033.467
$ws.Range("I132").Value = 5312.75
$ws.Range("K132").Value = 47814.75
$ws.Range("M132").Value = -45284.75
$ws.Range("H134").Value = 1531.7142
$ws.Range("I134").Value = 1531.7142
$ws.Range("K134").Value = 4595.142599999999
$ws.Range("M134").Value = 474.8574000000008
$ws.Range("H140").Value = 2742.923
$ws.Range("I140").Value = 2196.6365
$ws.Range("J140").Value = 5747.5
$ws.Range("K140").Value = 6589.9095
$ws.Range("L140").Value = 17242.5
$ws.Range("M140").Value = -1409.9095
$ws.Range("N140").Value = -27602.5
$ws.Range("H141").Value = 5149.5
$ws.Range("I141").Value = 4456.5713
$ws.Range("J141").Value = 10000
$ws.Range("K141").Value = 13369.7139
$ws.Range("L141").Value = 30000
$ws.Range("M141").Value = -8189.713899999999
$ws.Range("N141").Value = -40360

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 93958
$ws.Range("I132").Value = 93958
$ws.Range("K132").Value = 281874
$ws.Range("M132").Value = -279344

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 29999
$ws.Range("I43").Value = 29999
$ws.Range("K43").Value = 29999
$ws.Range("M43").Value = -29806
$ws.Range("H46").Value = 4571.5713
$ws.Range("J46").Value = 8000.3335
$ws.Range("L46").Value = 8000.3335
$ws.Range("N46").Value = -8376.333500000001
$ws.Range("H82").Value = 3372.182
$ws.Range("I82").Value = 1199.8
$ws.Range("J82").Value = 5182.5
$ws.Range("K82").Value = 1199.8
$ws.Range("L82").Value = 5182.5
$ws.Range("M82").Value = -838.8
$ws.Range("N82").Value = -5904.5
$ws.Range("H85").Value = 3372.182
$ws.Range("I85").Value = 1199.8
$ws.Range("J85").Value = 5182.5
$ws.Range("K85").Value = 1199.8
$ws.Range("L85").Value = 5182.5
$ws.Range("M85").Value = 48.20000000000005
$ws.Range("N85").Value = -7678.5
$ws.Range("H122").Value = 7076.154
$ws.Range("I122").Value = 6180.5454
$ws.Range("K122").Value = 18541.6362
$ws.Range("M122").Value = -16091.6362
$ws.Range("H132").Value = 4088.0908
$ws.Range("I132").Value = 2395.7334
$ws.Range("J132").Value = 7714.5713
$ws.Range("K132").Value = 7187.2002
$ws.Range("L132").Value = 23143.7139
$ws.Range("M132").Value = -4657.2002
$ws.Range("N132").Value = -28203.7139
$ws.Range("H136").Value = 6581.2383
$ws.Range("I136").Value = 3958
$ws.Range("K136").Value = 11874
$ws.Range("M136").Value = -9324

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2295.9375
$ws.Range("I122").Value = 1795.5172
$ws.Range("J122").Value = 7133.3335
$ws.Range("K122").Value = 5386.5516
$ws.Range("L122").Value = 21400.0005
$ws.Range("M122").Value = -2936.5516
$ws.Range("N122").Value = -26300.0005
$ws.Range("H126").Value = 6534.5884
$ws.Range("I126").Value = 5156.357
$ws.Range("J126").Value = 12966.333
$ws.Range("K126").Value = 15469.071
$ws.Range("L126").Value = 38898.999
$ws.Range("M126").Value = -12999.071
$ws.Range("N126").Value = -43838.999
$ws.Range("H132").Value = 5026.077
$ws.Range("I132").Value = 3585.125
$ws.Range("J132").Value = 7331.6
$ws.Range("K132").Value = 10755.375
$ws.Range("L132").Value = 21994.8
$ws.Range("M132").Value = -8225.375
$ws.Range("N132").Value = -27054.8
$ws.Range("H136").Value = 5414.6553
$ws.Range("I136").Value = 2903.1738
$ws.Range("J136").Value = 15042
$ws.Range("K136").Value = 8709.5214
$ws.Range("L136").Value = 45126
$ws.Range("M136").Value = -6159.5214
$ws.Range("N136").Value = -50226
